$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the end of the document (after
#    "Tracking de los resultados.") to straddle "Labview KEIVFRONTPANEL"
#    in the "Control de temperatura" paragraph. Adding a bookmark with
#    the same name elsewhere removes the old one (only one bookmark per
#    name can exist), and it also splits the run that used to read
#    " KEIVFRONTPANEL SIC (situado en la carpeta "2017_Control RD DAQ
#    v2")" into " KEIVFRONTPANEL" / " SIC (situado en la carpeta ...)".
# ------------------------------------------------------------------
$rngStart = $d.Content
$rngStart.Find.Execute("Labview", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$bmStart = $rngStart.Start

$rngEnd = $d.Content
$rngEnd.Find.Execute("KEIVFRONTPANEL", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$bmEnd = $rngEnd.End

$bmRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 2) Merge runs that only differ by an (insignificant) run split, by
#    replacing each span of text with itself -- Find/Replace collapses
#    the matched span into a single run.
# ------------------------------------------------------------------
$t1 = ".  En este caso, se sugiere usar valores de campo intercalados por medidas de ausencia de campo, es decir, B = [0 100 0 200 0 300 …] "
$d.Content.Find.Execute($t1, $true, $false, $false, $false, $false, $true, 1, $false, $t1, 2)

$t2 = " revisar la forma en la que se realizan los barridos en cada uno de los puntos anteriores, para poder minimizar el tiempo de caracterización. "
$d.Content.Find.Execute($t2, $true, $false, $false, $false, $false, $true, 1, $false, $t2, 2)

$t3 = " resolución adaptativa a los cambios que se producen en la curva I-V. La resolución debería de ser mayor en la zona de transición que fuera de ella para conseguir una mejor caracterización del fenómeno. En estos momentos, la resolución se fija manualmente pudiendo tener regiones con diferente resolución. Éste procedimiento requiere de la adquisición previa de una curva I-V con una resolución gruesa, para posteriormente optimizarla. "
$d.Content.Find.Execute($t3, $true, $false, $false, $false, $false, $true, 1, $false, $t3, 2)

# ------------------------------------------------------------------
# 3) Remove the two empty paragraphs right after "... Rn (función
#    RnCalc)." and before the final "IVset: ..." paragraph. Locate the
#    anchor paragraph by content (robust to any index shift) and then
#    delete the following two (empty) paragraph marks one at a time.
# ------------------------------------------------------------------
$n = $d.Paragraphs.Count
$anchorIdx = -1
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*RnCalc*") {
        $anchorIdx = $i
    }
}

$delIdx = $anchorIdx + 1
$pDel1 = $d.Paragraphs.Item($delIdx)
$pDel1.Range.Delete()
$pDel2 = $d.Paragraphs.Item($delIdx)
$pDel2.Range.Delete()
